# Refactor user maintenance component to include user role in table display
# Adds a "role" column to the users table (header + value for the sample row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F1").Value = "role"
$ws.Range("F2").Value = "admin"
